# "Add files via upload" — replace the placeholder industry-ticker row
# (FTSE, SS, EE, HFJF, DDS, ADA, WA, DA, EFFE, WW, AAD) with real LSE
# ticker symbols, matching the author's final worksheet state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds one ticker per industry column (row 1 = industry headers).
# Values are entered in the same order the author typed them (this is
# also the order their text first appears in the saved shared-string
# table): A, C, B, D, E, F, H, I, J, K, then G last.
$ws.Cells.Item(2, 1).Value  = "BREE.L"   # A2  Basic Materials
$ws.Cells.Item(2, 3).Value  = "CPG.L"    # C2  Consumer Cyclical
$ws.Cells.Item(2, 2).Value  = "VOD.L"    # B2  Communication Services
$ws.Cells.Item(2, 4).Value  = "DGE.L"    # D2  Consumer Defensive
$ws.Cells.Item(2, 5).Value  = "DCC.L"    # E2  Energy
$ws.Cells.Item(2, 6).Value  = "BBGI.L"   # F2  Financial Services
$ws.Cells.Item(2, 8).Value  = "REL.L"    # H2  Industrials
$ws.Cells.Item(2, 9).Value  = "BCPT.L"   # I2  Real Estate
$ws.Cells.Item(2, 10).Value = "SGE.L"    # J2  Technology
$ws.Cells.Item(2, 11).Value = "UU.L"     # K2  Utilities
$ws.Cells.Item(2, 7).Value  = "INDV.L"   # G2  Health Care

# Row 2 picked up a slightly custom height along the way.
$ws.Rows.Item(2).RowHeight = 13.8

# Page setup was touched (paper size / orientation explicitly set).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Final cursor position left on G3 (just below the last cell typed, G2).
[void]$ws.Range("G3").Select()
